$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Label" column header in H1, matching the style of the
#     other header cells (bold/centered/bordered style used by B1:G1) ---
$ws.Range("H1").Value() = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Refit values for the "100 iterations" block (rows 2-11): update the
#     Prediction/Error columns (D/E) and the one Cross Entropy Loss value
#     in F11 to the newly refit numbers ---
$ws.Range("D2").Value() = 0.2881363903532148
$ws.Range("E2").Value() = 0.2881363903532148

$ws.Range("D3").Value() = 0.4833938019656307
$ws.Range("E3").Value() = 0.4833938019656307

$ws.Range("D4").Value() = 0.4541178877172669
$ws.Range("E4").Value() = 0.4541178877172669

$ws.Range("D5").Value() = [double]"3.441234887121282E-14"
$ws.Range("E5").Value() = [double]"3.441234887121282E-14"

$ws.Range("D6").Value() = 0.1426756862535494
$ws.Range("E6").Value() = 0.1426756862535494

$ws.Range("D8").Value() = 0.4535478199234987
$ws.Range("E8").Value() = 0.5464521800765012

$ws.Range("D9").Value() = 0.4504288478113062
$ws.Range("E9").Value() = 0.5495711521886939

$ws.Range("D10").Value() = 0.4727668935396369
$ws.Range("E10").Value() = 0.5272331064603631

$ws.Range("D11").Value() = 0.000995361048090965
$ws.Range("E11").Value() = 0.9990046389519091
$ws.Range("F11").Value() = 1.138649106025696

# --- Populate the new "Label" column (H) for every data row. Rows 2-6 and
#     12-16 (Control patients) are labelled 0; rows 7-11 and 17-21 (MDD
#     patients) are labelled 1 ---
$ws.Range("H2").Value() = 0
$ws.Range("H3").Value() = 0
$ws.Range("H4").Value() = 0
$ws.Range("H5").Value() = 0
$ws.Range("H6").Value() = 0
$ws.Range("H7").Value() = 1
$ws.Range("H8").Value() = 1
$ws.Range("H9").Value() = 1
$ws.Range("H10").Value() = 1
$ws.Range("H11").Value() = 1
$ws.Range("H12").Value() = 0
$ws.Range("H13").Value() = 0
$ws.Range("H14").Value() = 0
$ws.Range("H15").Value() = 0
$ws.Range("H16").Value() = 0
$ws.Range("H17").Value() = 1
$ws.Range("H18").Value() = 1
$ws.Range("H19").Value() = 1
$ws.Range("H20").Value() = 1
$ws.Range("H21").Value() = 1
